# Generate Report for Handback
# Updates the handback-status workbook with new generated file names
# (uuid1 and uuid2 pseudo-GUIDs) and new handback timestamps, for both
# the zh-cn and de-de locales, across the Overview / zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

# --- new identifiers used to build the new file names -------------------
$uuid1New = 'd90d6a7c-5ffc-4db8-9e07-497e5621c334'
$uuid2New = 'ffffb0782bef-8311-4f94-aaa6-686e82cc04b3'
$hashZh   = '8c699e8997456c7ca9e599d4d73f26a578a91859'
$hashDe   = '8c699e8997456c7ca9e599d4d73f26a578a91859'

$md1New = "$uuid1New.md"
$md2New = "$uuid2New.md"

$xlfZhNew = "$uuid1New.$hashZh.zh-cn.xlf"
$xlfDeNew = "$uuid1New.$hashDe.de-de.xlf"

$zhDate1 = '2016-03-11 12:41:52'
$zhDate2 = '2016-03-11 12:42:10'
$deDate1 = '2016-03-11 12:41:56'
$deDate2 = '2016-03-11 12:42:15'

function Update-SheetCells {
    param($ws, $map)

    # 1) write the new cell values
    foreach ($addr in $map.Keys) {
        $plain = $addr.Replace('$', '')
        $ws.Range($plain).Value = $map[$addr]
    }

    # 2) keep hyperlink display text in sync with the new cell values,
    #    while preserving the existing hyperlink (address / r:id)
    foreach ($hl in $ws.Hyperlinks) {
        $rngAddr = $hl.Range.Address()
        if ($map.ContainsKey($rngAddr)) {
            $hl.TextToDisplay = $map[$rngAddr]
        }
    }
}

# --- Sheet "Overview" -----------------------------------------------------
$wsOverview = $wb.Worksheets.Item('Overview')
$overviewMap = @{
    '$A$2' = $md1New
    '$A$3' = $md2New
}
Update-SheetCells $wsOverview $overviewMap

# --- Sheet "zh-cn" ----------------------------------------------------------
$wsZh = $wb.Worksheets.Item('zh-cn')
$zhMap = @{
    '$A$2' = $md1New
    '$D$2' = $xlfZhNew
    '$E$2' = $zhDate1
    '$F$2' = $md1New
    '$G$2' = $xlfZhNew
    '$H$2' = $zhDate2
    '$A$3' = $md2New
    '$D$3' = $xlfZhNew
    '$E$3' = $zhDate1
    '$F$3' = $md2New
    '$G$3' = $xlfZhNew
    '$H$3' = $zhDate2
}
Update-SheetCells $wsZh $zhMap

# --- Sheet "de-de" ----------------------------------------------------------
$wsDe = $wb.Worksheets.Item('de-de')
$deMap = @{
    '$A$2' = $md1New
    '$D$2' = $xlfDeNew
    '$E$2' = $deDate1
    '$F$2' = $md1New
    '$G$2' = $xlfDeNew
    '$H$2' = $deDate2
    '$A$3' = $md2New
    '$D$3' = $xlfDeNew
    '$E$3' = $deDate1
    '$F$3' = $md2New
    '$G$3' = $xlfDeNew
    '$H$3' = $deDate2
}
Update-SheetCells $wsDe $deMap
